{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Net content change (per the supplied diff):\n//   1. The paragraph \"<Mitigation Type>\" (directly above the \"<Mitigation>\"\n//      paragraph, near the end of the letter) is removed entirely.\n//   2. The hidden \"_GoBack\" bookmark, which used to sit between the runs\n//      \"condition,\" and \" please let us know...\" earlier in the same\n//      section, is moved so that it now sits between the \"<\" and\n//      \"Mitigation\" runs of the (now final) \"<Mitigation>\" paragraph.\n\nconst body = context.document.body;\n\n// Step 1: remove the existing \"_GoBack\" bookmark from its old location.\n// deleteBookmark is a no-op/throws if the bookmark does not exist, so we\n// guard with getBookmarkRangeOrNullObject first for robustness.\nconst oldBmRange = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\nif (!oldBmRange.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// Step 2: find + delete the \"<Mitigation Type>\" paragraph.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet mitigationTypeParaIndex = -1;\nlet mitigationParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"<Mitigation Type>\") {\n    mitigationTypeParaIndex = i;\n  } else if (t === \"<Mitigation>\") {\n    mitigationParaIndex = i;\n  }\n}\n\nif (mitigationTypeParaIndex !== -1) {\n  paragraphs.items[mitigationTypeParaIndex].delete();\n  await context.sync();\n}\n\n// Step 3: re-locate the \"<Mitigation>\" paragraph (its index shifts down by\n// one once the paragraph above it was removed) and insert the \"_GoBack\"\n// bookmark right after its leading \"<\" run.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nlet finalMitigationIndex = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text === \"<Mitigation>\") {\n    finalMitigationIndex = i;\n    break;\n  }\n}\n\nif (finalMitigationIndex !== -1) {\n  const targetPara = paragraphs2.items[finalMitigationIndex];\n  const searchResults = targetPara.getRange().search(\"<\", { matchCase: true });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length > 0) {\n    const afterBracket = searchResults.items[0].getRange(\"After\");\n    afterBracket.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Net content change (per the supplied diff):\n#   1. The paragraph \"<Mitigation Type>\" (directly above the \"<Mitigation>\"\n#      paragraph, near the end of the letter) is removed entirely.\n#   2. The hidden \"_GoBack\" bookmark, which used to sit between the runs\n#      \"condition,\" and \" please let us know...\" earlier in the same\n#      section, is moved so that it now sits between the \"<\" and\n#      \"Mitigation\" runs of the (now final) \"<Mitigation>\" paragraph.\n\n$d = $word.ActiveDocument\n\n# Step 1: remove the \"_GoBack\" bookmark from its old location (if present).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Step 2: find and delete the \"<Mitigation Type>\" paragraph.\n$mitigationTypePara = $null\n$mitigationPara = $null\nforeach ($p in $d.Paragraphs) {\n  $text = $p.Range.Text\n  $trimmed = $text.TrimEnd([char]13, [char]7)\n  if ($trimmed -eq \"<Mitigation Type>\") {\n    $mitigationTypePara = $p\n  } elseif ($trimmed -eq \"<Mitigation>\") {\n    $mitigationPara = $p\n  }\n}\n\nif ($mitigationTypePara -ne $null) {\n  $mitigationTypePara.Range.Delete()\n}\n\n# Step 3: re-locate the \"<Mitigation>\" paragraph (it may have been\n# invalidated/shifted by the deletion above) and insert the \"_GoBack\"\n# bookmark right after its leading \"<\" character.\n$finalPara = $null\nforeach ($p in $d.Paragraphs) {\n  $text = $p.Range.Text\n  $trimmed = $text.TrimEnd([char]13, [char]7)\n  if ($trimmed -eq \"<Mitigation>\") {\n    $finalPara = $p\n    break\n  }\n}\n\nif ($finalPara -ne $null) {\n  $insertPoint = $finalPara.Range.Start + 1\n  $bmRange = $d.Range($insertPoint, $insertPoint)\n  $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
